# This script applies the changes described by the commit:
# "stop create report to update code for caculation lasary"
#
# Sheet "Đơn sale chính" (sheet 1):
#   - Insert a new column G "Nhóm dịch vụ" (shifting old G:V to H:W)
#   - Append 4 new trailing columns: X "Tỉ lệ chiết khấu sale chính",
#     Y "Tỉ lệ chiết khấu sale phụ", Z "Chiết khấu sale chính",
#     AA "Chiết khấu sale phụ"
#   - Fill in the new values for row 2 (data row) and row 3 (totals row)
#
# Sheet "Lương" (sheet 3):
#   - Populate a brand-new lookup table of salary/commission parameters
#     in columns A (label) and B (value), rows 1-24

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Đơn sale chính"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert the new "Nhóm dịch vụ" column at G; this shifts the existing
# G:V columns (and all their row 1-3 data) one slot to the right (H:W).
$ws1.Range("G1:G3").EntireColumn.Insert()

# --- Header row (row 1) ---
$ws1.Range("G1").Value = "Nhóm dịch vụ"
$ws1.Range("X1").Value = "Tỉ lệ chiết khấu sale chính"
$ws1.Range("Y1").Value = "Tỉ lệ chiết khấu sale phụ"
$ws1.Range("Z1").Value = "Chiết khấu sale chính"
$ws1.Range("AA1").Value = "Chiết khấu sale phụ"

# --- Data row (row 2) ---
$ws1.Range("G2").Value = "Vùng mắt"
# These four cells were shifted in from previously-empty numeric cells;
# the insert operation materialised them as 0, so restore true emptiness.
$ws1.Range("K2").ClearContents()
$ws1.Range("L2").ClearContents()
$ws1.Range("S2").ClearContents()
$ws1.Range("U2").ClearContents()
$ws1.Range("V2").Value = 50000
$ws1.Range("W2").Value = 0
$ws1.Range("X2").Value = 0.15
$ws1.Range("Y2").Value = 0
$ws1.Range("Z2").Value = 600000
$ws1.Range("AA2").Value = 0

# --- Totals row (row 3) ---
$ws1.Range("V3").Value = 50000
$ws1.Range("W3").Value = 0
$ws1.Range("X3").Value = 0.15
$ws1.Range("Y3").Value = 0
$ws1.Range("Z3").Value = 600000
$ws1.Range("AA3").Value = 0

# ---------------------------------------------------------------------
# Sheet 3: "Lương"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A1").Value = "Danh mục"
$ws3.Range("B1").Value = 11
$ws3.Range("A2").Value = "Ngày công"
$ws3.Range("B2").Value = 6
$ws3.Range("A3").Value = "Phụ cấp"
$ws3.Range("B3").Value = 210000
$ws3.Range("A4").Value = "Lương cơ bản tại CẦN THƠ"
$ws3.Range("A5").Value = "Chiết khấu sale chính tại CẦN THƠ"
$ws3.Range("B5").Value = 0
$ws3.Range("A6").Value = "Chiết khấu sale phụ tại CẦN THƠ"
$ws3.Range("B6").Value = 0
$ws3.Range("A7").Value = "Đơn 1 bác sĩ tại CẦN THƠ"
$ws3.Range("B7").Value = 0
$ws3.Range("A8").Value = "Đơn 2 bác sĩ tại CẦN THƠ"
$ws3.Range("B8").Value = 0
$ws3.Range("A9").Value = "Công phụ phẫu 1 tại CẦN THƠ"
$ws3.Range("B9").Value = 0
$ws3.Range("A10").Value = "Công phụ phẫu 2 tại CẦN THƠ"
$ws3.Range("B10").Value = 0
$ws3.Range("A11").Value = "Lương cơ bản tại LONG XUYÊN"
$ws3.Range("A12").Value = "Chiết khấu sale chính tại LONG XUYÊN"
$ws3.Range("B12").Value = 0
$ws3.Range("A13").Value = "Chiết khấu sale phụ tại LONG XUYÊN"
$ws3.Range("B13").Value = 0
$ws3.Range("A14").Value = "Đơn 1 bác sĩ tại LONG XUYÊN"
$ws3.Range("B14").Value = 0
$ws3.Range("A15").Value = "Đơn 2 bác sĩ tại LONG XUYÊN"
$ws3.Range("B15").Value = 0
$ws3.Range("A16").Value = "Công phụ phẫu 1 tại LONG XUYÊN"
$ws3.Range("B16").Value = 0
$ws3.Range("A17").Value = "Công phụ phẫu 2 tại LONG XUYÊN"
$ws3.Range("B17").Value = 0
$ws3.Range("A18").Value = "Lương cơ bản tại SÓC TRĂNG"
$ws3.Range("A19").Value = "Chiết khấu sale chính tại SÓC TRĂNG"
$ws3.Range("B19").Value = 600000
$ws3.Range("A20").Value = "Chiết khấu sale phụ tại SÓC TRĂNG"
$ws3.Range("B20").Value = 0
$ws3.Range("A21").Value = "Đơn 1 bác sĩ tại SÓC TRĂNG"
$ws3.Range("B21").Value = 0
$ws3.Range("A22").Value = "Đơn 2 bác sĩ tại SÓC TRĂNG"
$ws3.Range("B22").Value = 0
$ws3.Range("A23").Value = "Công phụ phẫu 1 tại SÓC TRĂNG"
$ws3.Range("B23").Value = 0
$ws3.Range("A24").Value = "Công phụ phẫu 2 tại SÓC TRĂNG"
$ws3.Range("B24").Value = 0

Write-Host "Edit applied: sheet1 columns inserted/appended, sheet3 populated."
